$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row Price (D) / Volume (E) updates ---
$ws.Range("D2").Value = '58.167.69'
$ws.Range("E2").Value = '  -3.54%  '

$ws.Range("D3").Value = '3.130.73'
$ws.Range("E3").Value = '  -5.18%  '

$ws.Range("D5").Value = '524.27'
$ws.Range("E5").Value = '  -5.89%  '

$ws.Range("D6").Value = '132.85'
$ws.Range("E6").Value = '  -5.98%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.129.35'
$ws.Range("E8").Value = '  -5.26%  '

$ws.Range("E9").Value = '  -4.77%  '

$ws.Range("D10").Value = '7.28'
$ws.Range("E10").Value = '  -8.10%  '

$ws.Range("E11").Value = '  -8.85%  '

$ws.Range("E12").Value = '  -8.61%  '

$ws.Range("D13").Value = '3.668.97'
$ws.Range("E13").Value = '  -5.08%  '

$ws.Range("E14").Value = '  -0.51%  '

$ws.Range("D15").Value = '25.02'
$ws.Range("E15").Value = '  -6.17%  '

$ws.Range("D16").Value = '3.134.58'
$ws.Range("E16").Value = '  -5.10%  '

$ws.Range("D17").Value = '58.181.47'
$ws.Range("E17").Value = '  -3.50%  '

$ws.Range("E18").Value = '  -7.95%  '

$ws.Range("D19").Value = '5.70'
$ws.Range("E19").Value = '  -5.82%  '

$ws.Range("D20").Value = '12.94'
$ws.Range("E20").Value = '  -5.74%  '

$ws.Range("D21").Value = '7.86'
$ws.Range("E21").Value = '  -7.81%  '

$ws.Range("D22").Value = '342.29'
$ws.Range("E22").Value = '  -8.44%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("E24").Value = '  -4.52%  '

$ws.Range("D25").Value = '67.18'

$ws.Range("D26").Value = '3.267.66'
$ws.Range("E26").Value = '  -4.81%  '

$ws.Range("E27").Value = '  -2.73%  '

$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.10%  '

$ws.Range("D29").Value = '0.0₃0923'
$ws.Range("E29").Value = '  -9.70%  '

$ws.Range("D32").Value = '1.27'
$ws.Range("E32").Value = '  +1.30%  '

$ws.Range("E33").Value = '  -7.92%  '

$ws.Range("D34").Value = '6.84'
$ws.Range("E34").Value = '  -7.81%  '

$ws.Range("D35").Value = '21.37'
$ws.Range("E35").Value = '  -5.28%  '

$ws.Range("D36").Value = '159.66'
$ws.Range("E36").Value = '  -3.65%  '

$ws.Range("D37").Value = '4.80'
$ws.Range("E37").Value = '  -4.66%  '

$ws.Range("D38").Value = '6.19'
$ws.Range("E38").Value = '  -6.50%  '

$ws.Range("E39").Value = '  -10.19%  '

$ws.Range("E42").Value = '  -2.96%  '

$ws.Range("D43").Value = '23.52'
$ws.Range("E43").Value = '  -7.70%  '

$ws.Range("D44").Value = '0.687'
$ws.Range("E44").Value = '  -8.22%  '

$ws.Range("D45").Value = '1.07'
$ws.Range("E45").Value = '  -4.76%  '

$ws.Range("E46").Value = '  -4.23%  '

$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("D48").Value = '1.44'
$ws.Range("E48").Value = '  -8.24%  '

$ws.Range("D49").Value = '2.276.24'
$ws.Range("E49").Value = '  -2.05%  '

$ws.Range("D50").Value = '6.14'
$ws.Range("E50").Value = '  -3.17%  '

$ws.Range("D51").Value = '20.42'
$ws.Range("E51").Value = '  -5.42%  '

# --- Row 30/31 swap: RenderToken <-> USDe (with updated values) ---
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "6.72"
$ws.Range("E31").Value = "  -4.45%  "

# --- Row 40/41 swap: Hedera <-> RenzoRestakedETH (with updated values) ---
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "3.163.68"
$ws.Range("E40").Value = "  -5.03%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.0680"
$ws.Range("E41").Value = "  -5.84%  "
